# Nudge a handful of shapes on slides 2-4 (connector lines + callout
# rectangles around the diagram) to their new positions.
#
# PowerPoint COM reports/accepts Shape.Left/.Top in points (1 pt =
# 12700 EMU); only the axis/axes that actually move are touched below
# so untouched coordinates round-trip bit-for-bit.

$p = $ppt.ActivePresentation

function Move-ShapeById($slide, $id, $left, $top) {
    foreach ($shp in $slide.Shapes) {
        if ($shp.Id -eq $id) {
            if ($left -ne $null) { $shp.Left = $left }
            if ($top -ne $null) { $shp.Top = $top }
        }
    }
}

$s2 = $p.Slides.Item(2)
$s3 = $p.Slides.Item(3)
$s4 = $p.Slides.Item(4)

# Slide 2, "Прямая соединительная линия 43" - shift right only
Move-ShapeById $s2 44 578.9765930175782 $null

# Slide 3, "Прямоугольник 17" - shift right only
Move-ShapeById $s3 18 628.0 $null

# Slide 3, "Прямая соединительная линия 41" - move up-left
Move-ShapeById $s3 42 341.97661417322837 84.53401574803149

# Slide 4, "Прямоугольник 17" - shift right only
Move-ShapeById $s4 18 629.9822834645669 $null

# Slide 4, "Прямоугольник 25" - move right and slightly down
Move-ShapeById $s4 26 629.9822834645669 149.78929133858267

# Slide 4, "Прямая соединительная линия 41" - shift up only
Move-ShapeById $s4 42 $null 68.88236220472442
